$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the level designations in column A: rows 34-65, 67, and 108-147
# change from 1 -> 2 (per commit: "updated the level designations in the DOE files")
$ws.Range("A34:A65").Value = 2
$ws.Range("A67").Value = 2
$ws.Range("A108:A147").Value = 2

# Update the active selection to reflect where the edit was made
$ws.Range("A108:A147").Select()
